$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("ttTopLevel"); this correctly shifts
# the old Name/Value/Phrase columns (C:E) to D:F, and carries the column
# width/style metadata along (column C joins the B:D "16.77734375" group).
$ws.Columns.Item(3).Insert()

# Only rows 1, 2, 3 and 7 actually receive content in the new column; the
# blank placeholder cells Insert() created on the other data rows are not
# part of the result, so remove them again completely.
$ws.Range("C4:C6").Clear() | Out-Null
$ws.Range("C8:C9").Clear() | Out-Null

# Row 1: plain label.
$ws.Range("C1").Value = "ttTopLevel"

# Rows 2, 3 and 7: formulas that echo column A of the same row.
$ws.Range("C2").Formula = '=$A2'
$ws.Range("C3").Formula = '=$A3'
$ws.Range("C7").Formula = '=$A7'

# The old "is toegekend" text (now shifted into E3) is dropped; keep the
# cell (and its formatting) but clear its contents.
$ws.Range("E3").ClearContents()

$ws.Range("C7").Select() | Out-Null
